# Day5-second commit: add patient details for the appointment
# Populates First Name / Last Name / email columns next to the existing
# "number" column, turns the email into a mailto: hyperlink (which picks
# up Excel's built-in "Hyperlink" style), widens the new columns to fit
# their content, and leaves the selection where the data entry ended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Last Name"
$ws.Range("D1").Value = "email"

# New patient record on row 2 (row 1 = headers, existing number rows stay put)
$ws.Range("B2").Value = "abi"
$ws.Range("C2").Value = "s"
$ws.Range("D2").Value = "abi@gmail.com"

# Turn the email address into a real mailto hyperlink
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:abi@gmail.com")

# Autosize the new columns around their new content
$ws.Columns.Item(2).ColumnWidth = 13.54296875
$ws.Columns.Item(3).ColumnWidth = 14.1796875
$ws.Columns.Item(4).ColumnWidth = 15.08984375
$ws.Columns.Item(5).ColumnWidth = 17.1796875

# Leave the selection where data entry finished
$ws.Range("B3").Select() | Out-Null
